# Lab 3 Questions.docx -- "finished question 4, can display post data"
#
# Inserts a new bold answer paragraph for Question 4 ("How does the
# POSTed data come to the CGI script?"), placed between the two existing
# blank paragraphs that already separate Question 4 from Question 5.

$d = $word.ActiveDocument

# Locate Question 4's text so the insertion point is found robustly
# instead of relying on a hard-coded paragraph index.
$q4 = $d.Content
$found = $q4.Find.Execute("data come to the CGI script?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find Question 4 text to anchor the new answer paragraph."
}

# Collapse to just after the found text, then move past the first blank
# paragraph that follows Question 4 -- that's where the new answer
# paragraph needs to be inserted (keeping both existing blank paragraphs
# intact, one before and one after the new answer).
$q4.Collapse(0)
$afterFirstBlank = $q4.Next(4, 1)
$insertPoint = $afterFirstBlank.End

# Create a brand-new (still empty) paragraph right there.
$null = $afterFirstBlank.InsertParagraphAfter()

# Grab exactly that new paragraph (the mark that now lives at
# [$insertPoint, $insertPoint+1)) and fill it in with the bold answer
# text, run by run, matching the authored formatting.
$newParaRange = $d.Range($insertPoint, $insertPoint + 1)

$answerXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">The ' + [char]0x201C + 'FieldStorage' + [char]0x201D + ' class stores the </w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>data filled out in the</w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> username &amp; password </w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">fields </w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>and then retrieves t</w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>hose</w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> values using the method ' + [char]0x201C + '.getvalue()' + [char]0x201D + ' </w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>in order to display the username and password posted on the CGI script.</w:t></w:r>' +
  '</w:p>'

$null = $newParaRange.InsertXML($answerXml)

Write-Output "Inserted Question 4 answer paragraph."
